$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Itgav"
$ws.Range("C2").Value = "Thy1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 15.35884066666667
$ws.Range("H2").Value = 46.076522
$ws.Range("I2").Value = 0.1012042817263867
$ws.Range("J2").Value = 0.1012042817263867
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.114886
$ws.Range("N2").Value = 3.344658
$ws.Range("O2").Value = 0.01026455462504307
$ws.Range("P2").Value = 0.01026455462504307
$ws.Range("Q2").Value = 17.12335643549733
$ws.Range("R2").Value = 154.110207919476
$ws.Range("S2").Value = 0.001038816878068744
$ws.Range("T2").Value = 0.001038816878068744

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Itgav"
$ws.Range("C3").Value = "Thy1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.35884066666667
$ws.Range("H3").Value = 46.076522
$ws.Range("I3").Value = 0.1012042817263867
$ws.Range("J3").Value = 0.1012042817263867
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 80.666326
$ws.Range("N3").Value = 241.998978
$ws.Range("O3").Value = 0.7426803364904859
$ws.Range("P3").Value = 0.7426803364904858
$ws.Range("Q3").Value = 1238.941248199391
$ws.Range("R3").Value = 11150.47123379451
$ws.Range("S3").Value = 0.07516243000683077
$ws.Range("T3").Value = 0.07516243000683077

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Itgav"
$ws.Range("C4").Value = "Thy1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 15.35884066666667
$ws.Range("H4").Value = 46.076522
$ws.Range("I4").Value = 0.1012042817263867
$ws.Range("J4").Value = 0.1012042817263867
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.6900326666666666
$ws.Range("N4").Value = 2.070098
$ws.Range("O4").Value = 0.006353006495788928
$ws.Range("P4").Value = 0.006353006495788927
$ws.Range("Q4").Value = 10.59810178212844
$ws.Range("R4").Value = 95.38291603915599
$ws.Range("S4").Value = 0.0006429514592093871
$ws.Range("T4").Value = 0.0006429514592093871

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Itgav"
$ws.Range("C5").Value = "Thy1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.35884066666667
$ws.Range("H5").Value = 46.076522
$ws.Range("I5").Value = 0.1012042817263867
$ws.Range("J5").Value = 0.1012042817263867
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.14389166666666
$ws.Range("N5").Value = 78.431675
$ws.Range("O5").Value = 0.2407021023886821
$ws.Range("P5").Value = 0.2407021023886821
$ws.Range("Q5").Value = 401.5398665149277
$ws.Range("R5").Value = 3613.85879863435
$ws.Range("S5").Value = 0.02436008338227775
$ws.Range("T5").Value = 0.02436008338227775

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Itgav"
$ws.Range("C6").Value = "Thy1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 50.59256466666667
$ws.Range("H6").Value = 151.777694
$ws.Range("I6").Value = 0.3333704853712116
$ws.Range("J6").Value = 0.3333704853712116
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.114886
$ws.Range("N6").Value = 3.344658
$ws.Range("O6").Value = 0.01026455462504307
$ws.Range("P6").Value = 0.01026455462504307
$ws.Range("Q6").Value = 56.40494205096134
$ws.Range("R6").Value = 507.644478458652
$ws.Range("S6").Value = 0.003421899557469922
$ws.Range("T6").Value = 0.003421899557469921

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Itgav"
$ws.Range("C7").Value = "Thy1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 50.59256466666667
$ws.Range("H7").Value = 151.777694
$ws.Range("I7").Value = 0.3333704853712116
$ws.Range("J7").Value = 0.3333704853712116
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 80.666326
$ws.Range("N7").Value = 241.998978
$ws.Range("O7").Value = 0.7426803364904859
$ws.Range("P7").Value = 0.7426803364904858
$ws.Range("Q7").Value = 4081.116314577415
$ws.Range("R7").Value = 36730.04683119673
$ws.Range("S7").Value = 0.247587704251488
$ws.Range("T7").Value = 0.247587704251488

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Itgav"
$ws.Range("C8").Value = "Thy1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 50.59256466666667
$ws.Range("H8").Value = 151.777694
$ws.Range("I8").Value = 0.3333704853712116
$ws.Range("J8").Value = 0.3333704853712116
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.6900326666666666
$ws.Range("N8").Value = 2.070098
$ws.Range("O8").Value = 0.006353006495788928
$ws.Range("P8").Value = 0.006353006495788927
$ws.Range("Q8").Value = 34.91052231044578
$ws.Range("R8").Value = 314.194700794012
$ws.Range("S8").Value = 0.002117904859067615
$ws.Range("T8").Value = 0.002117904859067614

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Itgav"
$ws.Range("C9").Value = "Thy1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 50.59256466666667
$ws.Range("H9").Value = 151.777694
$ws.Range("I9").Value = 0.3333704853712116
$ws.Range("J9").Value = 0.3333704853712116
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 26.14389166666666
$ws.Range("N9").Value = 78.431675
$ws.Range("O9").Value = 0.2407021023886821
$ws.Range("P9").Value = 0.2407021023886821
$ws.Range("Q9").Value = 1322.686529784161
$ws.Range("R9").Value = 11904.17876805745
$ws.Range("S9").Value = 0.08024297670318602
$ws.Range("T9").Value = 0.080242976703186

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Itgav"
$ws.Range("C10").Value = "Thy1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 60.37715666666667
$ws.Range("H10").Value = 181.13147
$ws.Range("I10").Value = 0.397844271305776
$ws.Range("J10").Value = 0.397844271305776
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.114886
$ws.Range("N10").Value = 3.344658
$ws.Range("O10").Value = 0.01026455462504307
$ws.Range("P10").Value = 0.01026455462504307
$ws.Range("Q10").Value = 67.31364668747334
$ws.Range("R10").Value = 605.82282018726
$ws.Range("S10").Value = 0.004083694255078592
$ws.Range("T10").Value = 0.004083694255078592

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Itgav"
$ws.Range("C11").Value = "Thy1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 60.37715666666667
$ws.Range("H11").Value = 181.13147
$ws.Range("I11").Value = 0.397844271305776
$ws.Range("J11").Value = 0.397844271305776
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 80.666326
$ws.Range("N11").Value = 241.998978
$ws.Range("O11").Value = 0.7426803364904859
$ws.Range("P11").Value = 0.7426803364904858
$ws.Range("Q11").Value = 4870.403402626407
$ws.Range("R11").Value = 43833.63062363766
$ws.Range("S11").Value = 0.2954711172841859
$ws.Range("T11").Value = 0.2954711172841858

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Itgav"
$ws.Range("C12").Value = "Thy1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 60.37715666666667
$ws.Range("H12").Value = 181.13147
$ws.Range("I12").Value = 0.397844271305776
$ws.Range("J12").Value = 0.397844271305776
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.6900326666666666
$ws.Range("N12").Value = 2.070098
$ws.Range("O12").Value = 0.006353006495788928
$ws.Range("P12").Value = 0.006353006495788927
$ws.Range("Q12").Value = 41.66221042045111
$ws.Range("R12").Value = 374.95989378406
$ws.Range("S12").Value = 0.002527507239918007
$ws.Range("T12").Value = 0.002527507239918007

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Itgav"
$ws.Range("C13").Value = "Thy1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 60.37715666666667
$ws.Range("H13").Value = 181.13147
$ws.Range("I13").Value = 0.397844271305776
$ws.Range("J13").Value = 0.397844271305776
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 26.14389166666666
$ws.Range("N13").Value = 78.431675
$ws.Range("O13").Value = 0.2407021023886821
$ws.Range("P13").Value = 0.2407021023886821
$ws.Range("Q13").Value = 1578.493843034694
$ws.Range("R13").Value = 14206.44458731225
$ws.Range("S13").Value = 0.09576195252659352
$ws.Range("T13").Value = 0.0957619525265935

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Itgav"
$ws.Range("C14").Value = "Thy1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 25.43221733333333
$ws.Range("H14").Value = 76.29665199999999
$ws.Range("I14").Value = 0.1675809615966257
$ws.Range("J14").Value = 0.1675809615966258
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 1.114886
$ws.Range("N14").Value = 3.344658
$ws.Range("O14").Value = 0.01026455462504307
$ws.Range("P14").Value = 0.01026455462504307
$ws.Range("Q14").Value = 28.35402305389066
$ws.Range("R14").Value = 255.186207485016
$ws.Range("S14").Value = 0.00172014393442581
$ws.Range("T14").Value = 0.00172014393442581

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Itgav"
$ws.Range("C15").Value = "Thy1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 25.43221733333333
$ws.Range("H15").Value = 76.29665199999999
$ws.Range("I15").Value = 0.1675809615966257
$ws.Range("J15").Value = 0.1675809615966258
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 80.666326
$ws.Range("N15").Value = 241.998978
$ws.Range("O15").Value = 0.7426803364904859
$ws.Range("P15").Value = 0.7426803364904858
$ws.Range("Q15").Value = 2051.523534313517
$ws.Range("R15").Value = 18463.71180882166
$ws.Range("S15").Value = 0.1244590849479812
$ws.Range("T15").Value = 0.1244590849479812

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Itgav"
$ws.Range("C16").Value = "Thy1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 25.43221733333333
$ws.Range("H16").Value = 76.29665199999999
$ws.Range("I16").Value = 0.1675809615966257
$ws.Range("J16").Value = 0.1675809615966258
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.6900326666666666
$ws.Range("N16").Value = 2.070098
$ws.Range("O16").Value = 0.006353006495788928
$ws.Range("P16").Value = 0.006353006495788927
$ws.Range("Q16").Value = 17.54906074576622
$ws.Range("R16").Value = 157.941546711896
$ws.Range("S16").Value = 0.001064642937593918
$ws.Range("T16").Value = 0.001064642937593918

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Itgav"
$ws.Range("C17").Value = "Thy1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 25.43221733333333
$ws.Range("H17").Value = 76.29665199999999
$ws.Range("I17").Value = 0.1675809615966257
$ws.Range("J17").Value = 0.1675809615966258
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 26.14389166666666
$ws.Range("N17").Value = 78.431675
$ws.Range("O17").Value = 0.2407021023886821
$ws.Range("P17").Value = 0.2407021023886821
$ws.Range("Q17").Value = 664.8971348057888
$ws.Range("R17").Value = 5984.074213252099
$ws.Range("S17").Value = 0.04033708977662482
$ws.Range("T17").Value = 0.04033708977662482
